# Add Denmark, Sweden and Norway market test data sheets by cloning the
# Belgium sheet (closest matching layout/style) and adjusting the cells
# that differ per-market.

$wb = $excel.ActiveWorkbook
$belgium = $wb.Worksheets.Item("Belgium")

# --- Denmark -----------------------------------------------------------
$belgium.Copy($null, $belgium)
$denmark = $wb.Worksheets.Item($belgium.Index() + 1)
$denmark.Name = "Denmark"
$denmark.Rows("9:10").Delete()
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").ClearContents()

# --- Sweden --------------------------------------------------------------
$belgium.Copy($null, $denmark)
$sweden = $wb.Worksheets.Item($denmark.Index() + 1)
$sweden.Name = "Sweden"
$sweden.Rows("9:10").Delete()
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").ClearContents()
$sweden.Range("A8").Value = "FC604S"
$sweden.Range("G8").Value = "5.000"
$sweden.Range("K8").Value = "5.000"
$sweden.Range("M8").Value = "3.000"
$sweden.Range("O8").Value = "5.000"
$sweden.Range("R8").Value = 2200

# --- Norway --------------------------------------------------------------
$belgium.Copy($null, $sweden)
$norway = $wb.Worksheets.Item($sweden.Index() + 1)
$norway.Name = "Norway"
$norway.Rows("9:10").Delete()
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").ClearContents()

# --- View / selection state ----------------------------------------------
$denmark.Activate()
$denmark.Cells.Select()

$sweden.Activate()
$sweden.Range("E8").Select()

$norway.Activate()
$norway.Range("B3").Select()

$wb.Save()
